$d = $word.ActiveDocument

# Locate the placeholder paragraph text (" To be filled") that needs to be
# replaced with the real milestone description.
$finder = $d.Content
$found = $finder.Find.Execute(" To be filled", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the placeholder text ' To be filled' in the document."
}

# Re-materialize a fresh Range over the exact span Find located; re-using
# the Find range object directly with InsertXML does not replace content
# in-place, so a new Range bound to the same start/end is used instead.
$rng = $d.Range($finder.Start, $finder.End)

# Replace the whole run's content with the new set of runs, preserving the
# leading space run and adding the proofing-error markers around
# "Milestone" exactly as produced by Word's grammar checker.
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">The smartwatch will do some basic functions within the </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>Milestone</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> and I am planning to extend it beyond the milestone as well.</w:t></w:r>' +
    '</w:p>'
$rng.InsertXML($newParaXml)

# Add a new, empty paragraph right after that paragraph (before the section
# properties), matching the trailing "<w:p/>" added in the target document.
$d.Paragraphs.Add() | Out-Null
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>')
